$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is the target cell reference and its new text value, taken from
# the authoritative diff of the refreshed cryptocurrency listing.
$updates = @(
    @{ Ref = 'D2'; Val = '34.665.13' },
    @{ Ref = 'E2'; Val = '  +1.71%  ' },
    @{ Ref = 'D3'; Val = '1.807.66' },
    @{ Ref = 'D4'; Val = '1.00' },
    @{ Ref = 'E4'; Val = '  -0.02%  ' },
    @{ Ref = 'D5'; Val = '225.40' },
    @{ Ref = 'E5'; Val = '  -0.62%  ' },
    @{ Ref = 'D6'; Val = '0.554' },
    @{ Ref = 'E6'; Val = '  -0.30%  ' },
    @{ Ref = 'D7'; Val = '1.00' },
    @{ Ref = 'E7'; Val = '  -0.03%  ' },
    @{ Ref = 'E8'; Val = '  +4.75%  ' },
    @{ Ref = 'D9'; Val = '0.290' },
    @{ Ref = 'E9'; Val = '  +3.14%  ' },
    @{ Ref = 'E10'; Val = '  +8.02%  ' },
    @{ Ref = 'E11'; Val = '  +0.10%  ' },
    @{ Ref = 'D12'; Val = '2.067.79' },
    @{ Ref = 'E12'; Val = '  +1.07%  ' },
    @{ Ref = 'D13'; Val = '11.14' },
    @{ Ref = 'E13'; Val = '  -2.38%  ' },
    @{ Ref = 'D14'; Val = '1.814.98' },
    @{ Ref = 'E14'; Val = '  +1.41%  ' },
    @{ Ref = 'D15'; Val = '0.644' },
    @{ Ref = 'D16'; Val = '34.683.34' },
    @{ Ref = 'E16'; Val = '  +1.84%  ' },
    @{ Ref = 'D17'; Val = '4.33' },
    @{ Ref = 'E17'; Val = '  +2.55%  ' },
    @{ Ref = 'D18'; Val = '69.60' },
    @{ Ref = 'E18'; Val = '  +0.14%  ' },
    @{ Ref = 'D19'; Val = '254.34' },
    @{ Ref = 'E19'; Val = '  +0.38%  ' },
    @{ Ref = 'E20'; Val = '  +8.38%  ' },
    @{ Ref = 'D21'; Val = '10.95' },
    @{ Ref = 'E21'; Val = '  +4.75%  ' },
    @{ Ref = 'E22'; Val = '  -0.06%  ' },
    @{ Ref = 'E23'; Val = '  -0.83%  ' },
    @{ Ref = 'D25'; Val = '161.78' },
    @{ Ref = 'E25'; Val = '  +3.32%  ' },
    @{ Ref = 'D26'; Val = '16.50' },
    @{ Ref = 'E26'; Val = '  -0.61%  ' },
    @{ Ref = 'D27'; Val = '7.17' },
    @{ Ref = 'E27'; Val = '  +1.85%  ' },
    @{ Ref = 'D28'; Val = '0.114' },
    @{ Ref = 'E28'; Val = '  +0.20%  ' },
    @{ Ref = 'B29'; Val = 'Swop.fi' },
    @{ Ref = 'C29'; Val = 'https://coinranking.com/coin/yrCr2HW2c+swopfi-swop' },
    @{ Ref = 'D29'; Val = '649.23' },
    @{ Ref = 'E29'; Val = '  +1,129.58%  ' },
    @{ Ref = 'B30'; Val = 'BinanceUSD' },
    @{ Ref = 'C30'; Val = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd' },
    @{ Ref = 'D30'; Val = '0.999' },
    @{ Ref = 'E30'; Val = '  -0.07%  ' },
    @{ Ref = 'B31'; Val = 'Hedera' },
    @{ Ref = 'C31'; Val = 'https://coinranking.com/coin/jad286TjB+hedera-hbar' },
    @{ Ref = 'D31'; Val = '0.0534' },
    @{ Ref = 'E31'; Val = '  +3.31%  ' },
    @{ Ref = 'B32'; Val = 'Filecoin' },
    @{ Ref = 'C32'; Val = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil' },
    @{ Ref = 'D32'; Val = '3.81' },
    @{ Ref = 'E32'; Val = '  -0.47%  ' },
    @{ Ref = 'B33'; Val = 'PancakeSwap' },
    @{ Ref = 'C33'; Val = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake' },
    @{ Ref = 'D33'; Val = '1.21' },
    @{ Ref = 'E33'; Val = '  +0.20%  ' },
    @{ Ref = 'B34'; Val = 'InternetComputer(DFINITY)' },
    @{ Ref = 'C34'; Val = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp' },
    @{ Ref = 'D34'; Val = '3.64' },
    @{ Ref = 'E34'; Val = '  +0.77%  ' },
    @{ Ref = 'B35'; Val = 'LidoDAOToken' },
    @{ Ref = 'C35'; Val = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo' },
    @{ Ref = 'D35'; Val = '1.90' },
    @{ Ref = 'E35'; Val = '  +2.51%  ' },
    @{ Ref = 'B36'; Val = 'Maker' },
    @{ Ref = 'C36'; Val = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr' },
    @{ Ref = 'D36'; Val = '1.439.17' },
    @{ Ref = 'E36'; Val = '  -0.81%  ' },
    @{ Ref = 'B37'; Val = 'TrustWalletToken' },
    @{ Ref = 'C37'; Val = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt' },
    @{ Ref = 'D37'; Val = '1.07' },
    @{ Ref = 'E37'; Val = '  +0.38%  ' },
    @{ Ref = 'B39'; Val = 'ImmutableX' },
    @{ Ref = 'C39'; Val = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx' },
    @{ Ref = 'D39'; Val = '0.643' },
    @{ Ref = 'E39'; Val = '  +1.62%  ' },
    @{ Ref = 'B40'; Val = 'Aave' },
    @{ Ref = 'C40'; Val = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave' },
    @{ Ref = 'D40'; Val = '85.06' },
    @{ Ref = 'E40'; Val = '  +2.03%  ' },
    @{ Ref = 'B41'; Val = 'ARBITRUM' },
    @{ Ref = 'C41'; Val = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb' },
    @{ Ref = 'D41'; Val = '0.957' },
    @{ Ref = 'E41'; Val = '  +6.17%  ' },
    @{ Ref = 'B42'; Val = 'MXToken' },
    @{ Ref = 'C42'; Val = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' },
    @{ Ref = 'D42'; Val = '2.80' },
    @{ Ref = 'E42'; Val = '  -0.63%  ' },
    @{ Ref = 'B43'; Val = 'HuobiToken' },
    @{ Ref = 'C43'; Val = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht' },
    @{ Ref = 'D43'; Val = '2.34' },
    @{ Ref = 'E43'; Val = '  -0.17%  ' },
    @{ Ref = 'B44'; Val = 'RenderToken' },
    @{ Ref = 'C44'; Val = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' },
    @{ Ref = 'D44'; Val = '2.16' },
    @{ Ref = 'E44'; Val = '  +2.97%  ' },
    @{ Ref = 'B45'; Val = 'FraxShare' },
    @{ Ref = 'C45'; Val = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs' },
    @{ Ref = 'D45'; Val = '6.03' },
    @{ Ref = 'E45'; Val = '  +5.17%  ' },
    @{ Ref = 'B46'; Val = 'WEMIXToken' },
    @{ Ref = 'C46'; Val = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix' },
    @{ Ref = 'D46'; Val = '1.06' },
    @{ Ref = 'E46'; Val = '  -0.92%  ' },
    @{ Ref = 'B47'; Val = 'Kaspa' },
    @{ Ref = 'C47'; Val = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas' },
    @{ Ref = 'D47'; Val = '0.0494' },
    @{ Ref = 'E47'; Val = '  -3.46%  ' },
    @{ Ref = 'B48'; Val = 'RocketPoolETH' },
    @{ Ref = 'C48'; Val = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth' },
    @{ Ref = 'D48'; Val = '1.962.88' },
    @{ Ref = 'E48'; Val = '  +0.84%  ' },
    @{ Ref = 'B49'; Val = 'Quant' },
    @{ Ref = 'C49'; Val = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt' },
    @{ Ref = 'D49'; Val = '106.36' },
    @{ Ref = 'E49'; Val = '  +8.96%  ' },
    @{ Ref = 'B50'; Val = 'InjectiveProtocol' },
    @{ Ref = 'C50'; Val = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj' },
    @{ Ref = 'D50'; Val = '12.24' },
    @{ Ref = 'E50'; Val = '  +3.14%  ' },
    @{ Ref = 'B51'; Val = 'PaxDollar' },
    @{ Ref = 'C51'; Val = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp' },
    @{ Ref = 'D51'; Val = '1.00' },
    @{ Ref = 'E51'; Val = '  +0.08%  ' }
)

foreach ($u in $updates) {
    $ref = $u.Ref
    $val = $u.Val
    $range = $ws.Range($ref)

    # Columns B and C always hold plain text (coin names / URLs) so they can be
    # written as-is. Columns D/E hold price & percentage text that often LOOKS
    # like a number (e.g. "1.00", "225.40", "  +1.71%  "). Excel's COM layer
    # auto-coerces such strings to real numbers when assigned to a
    # General-formatted cell, which would silently destroy the intended
    # trailing-zero / thousands-dot text formatting. Detect that case and
    # prefix with an apostrophe (quote-prefix) to force text storage, exactly
    # as a human typing into Excel would need to.
    $looksNumeric = $val -match '^[+-]?[0-9]*\.?[0-9]+$'

    if ($looksNumeric) {
        $range.Value = "'" + $val
        # Re-normalize the style so the quote-prefix marker Excel just
        # attached doesn't leave a stray style index behind; the cell
        # keeps its text value but reverts to the sheet's default look,
        # matching how these cells were styled before the edit.
        $range.Style = "Normal"
    } else {
        $range.Value = $val
    }
}
